$p = $ppt.ActivePresentation

# Slide 2: TextBox 3 "The Moon" -> merge runs into a single run
$s2 = $p.Slides.Item(2)
$tb2 = $s2.Shapes.Item(2)
$tb2.TextFrame.TextRange.Text = "x"
$tb2.TextFrame.TextRange.Text = "The Moon"

# Slide 3: Title 1 "One More" -> merge runs into a single run
$s3 = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item(1)
$title3.TextFrame.TextRange.Text = "x"
$title3.TextFrame.TextRange.Text = "One More"

# Slide 3: TextBox 3 "The Moon" -> merge runs into a single run
$tb3 = $s3.Shapes.Item(3)
$tb3.TextFrame.TextRange.Text = "x"
$tb3.TextFrame.TextRange.Text = "The Moon"
